$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Re-arrange / rename worksheets so the tab order becomes:
#      总计, 2022-Q3, 2022-Q1, 2021-Q3, 2021-Q1
#    Rename from the back forward so we never collide with an existing name,
#    and insert the brand-new "2021-Q3" sheet right before "2021-Q1" (its
#    data used to live on the sheet that is now "2022-Q1").
# ---------------------------------------------------------------------------
$wsOldQ1_21   = $wb.Worksheets.Item("2021-Q1")   # keeps name 2021-Q1, stays last, untouched
$wsOldQ3_21   = $wb.Worksheets.Item("2021-Q3")   # will become "2022-Q1"
$wsOldQ1_22   = $wb.Worksheets.Item("2022-Q1")   # will become "2022-Q3"

# Remember the style donor BEFORE anything else changes (old "2022-Q1" sheet
# still has its original header/index-column formatting at this point).
$styleDonor = $wsOldQ1_22

$wsOldQ1_22.Name = "2022-Q3"
$wsOldQ3_21.Name = "2022-Q1"

$newQ3_21 = $wb.Worksheets.Add($null, $wsOldQ3_21)   # insert brand-new sheet right after "2022-Q1"
$newQ3_21.Name = "2021-Q3"

# ---------------------------------------------------------------------------
# 2. "总计" summary sheet
# ---------------------------------------------------------------------------
$zj = $wb.Worksheets.Item("总计")

# New row 5 needs the same index-column style ("s=2"-equivalent) as the
# existing rows; clone it from row 4 before overwriting values.
$zj.Range("A4").Copy($zj.Range("A5"))

$zj.Range("A2").Value = 0
$zj.Range("B2").Value = "2022-Q3"
$zj.Range("C2").Value = 2
$zj.Range("D2").Value = 0.07000000000000001

$zj.Range("A3").Value = 1
$zj.Range("B3").Value = "2022-Q1"
$zj.Range("C3").Value = 1
$zj.Range("D3").Value = 0.08

$zj.Range("A4").Value = 2
$zj.Range("B4").Value = "2021-Q3"
$zj.Range("C4").Value = 1
$zj.Range("D4").Value = 0.09

$zj.Range("A5").Value = 3
$zj.Range("B5").Value = "2021-Q1"
$zj.Range("C5").Value = 4
$zj.Range("D5").Value = 0.1

# ---------------------------------------------------------------------------
# 3. "2022-Q3" sheet (new quarter data) - this used to be the "2022-Q1" sheet
#    and already carries the correct header / index-column styling; clone
#    that same styling onto the brand-new row 3 before writing values.
# ---------------------------------------------------------------------------
$q3_22 = $wb.Worksheets.Item("2022-Q3")

$q3_22.Range("A2").Copy($q3_22.Range("A3"))
$q3_22.Range("B2:H2").Copy($q3_22.Range("B3:H3"))

$q3_22.Range("A2").Value = 0
$q3_22.Range("B2").NumberFormat = "@"
$q3_22.Range("B2").Value = "002295"
$q3_22.Range("C2").Value = "广发稳安灵活配置混合A"
$q3_22.Range("D2").NumberFormat = "@"
$q3_22.Range("D2").Value = "1.58"
$q3_22.Range("E2").NumberFormat = "@"
$q3_22.Range("E2").Value = "69.63"
$q3_22.Range("F2").NumberFormat = "@"
$q3_22.Range("F2").Value = "4.26"
$q3_22.Range("G2").NumberFormat = "@"
$q3_22.Range("G2").Value = "0.0673"
$q3_22.Range("H2").Value = 6

$q3_22.Range("A3").Value = 1
$q3_22.Range("B3").NumberFormat = "@"
$q3_22.Range("B3").Value = "008604"
$q3_22.Range("C3").Value = "广发稳安灵活配置混合C"
$q3_22.Range("D3").NumberFormat = "@"
$q3_22.Range("D3").Value = "0.02"
$q3_22.Range("E3").NumberFormat = "@"
$q3_22.Range("E3").Value = "69.63"
$q3_22.Range("F3").NumberFormat = "@"
$q3_22.Range("F3").Value = "4.26"
$q3_22.Range("G3").NumberFormat = "@"
$q3_22.Range("G3").Value = "0.0009"
$q3_22.Range("H3").Value = 6

# ---------------------------------------------------------------------------
# 4. "2022-Q1" sheet (previously held the 2021-Q3 fund, now holds the old
#    2022-Q1 fund, and its amount-column header becomes "基金规模"). It
#    already has the right formatting, so only values need writing.
# ---------------------------------------------------------------------------
$q1_22 = $wb.Worksheets.Item("2022-Q1")

$q1_22.Range("D1").Value = "基金规模"

$q1_22.Range("A2").Value = 0
$q1_22.Range("B2").NumberFormat = "@"
$q1_22.Range("B2").Value = "233009"
$q1_22.Range("C2").Value = "大摩多因子精选策略混合"
$q1_22.Range("D2").NumberFormat = "@"
$q1_22.Range("D2").Value = "6.77"
$q1_22.Range("E2").NumberFormat = "@"
$q1_22.Range("E2").Value = "89.73"
$q1_22.Range("F2").NumberFormat = "@"
$q1_22.Range("F2").Value = "1.18"
$q1_22.Range("G2").NumberFormat = "@"
$q1_22.Range("G2").Value = "0.0799"
$q1_22.Range("H2").Value = 5

# ---------------------------------------------------------------------------
# 5. "2021-Q3" sheet (brand new sheet; holds the old 2021-Q3 fund data that
#    used to live on the sheet now renamed to "2022-Q1"). Clone the header /
#    index-column styling from the style donor ("2022-Q3" sheet, formerly
#    "2022-Q1") before writing values.
# ---------------------------------------------------------------------------
$q3_21 = $wb.Worksheets.Item("2021-Q3")

$styleDonor.Range("B1:H1").Copy($q3_21.Range("B1:H1"))
$styleDonor.Range("A2").Copy($q3_21.Range("A2"))
$styleDonor.Range("B2:H2").Copy($q3_21.Range("B2:H2"))

$q3_21.Range("B1").Value = "基金代码"
$q3_21.Range("C1").Value = "基金名称"
$q3_21.Range("D1").Value = "基金金额"
$q3_21.Range("E1").Value = "股票总仓位"
$q3_21.Range("F1").Value = "仓位占比"
$q3_21.Range("G1").Value = "持有市值(亿元)"
$q3_21.Range("H1").Value = "仓位排名"

$q3_21.Range("A2").Value = 0
$q3_21.Range("B2").NumberFormat = "@"
$q3_21.Range("B2").Value = "004209"
$q3_21.Range("C2").Value = "大成智惠量化多策略灵活配置混合"
$q3_21.Range("D2").NumberFormat = "@"
$q3_21.Range("D2").Value = "1.40"
$q3_21.Range("E2").NumberFormat = "@"
$q3_21.Range("E2").Value = "92.89"
$q3_21.Range("F2").NumberFormat = "@"
$q3_21.Range("F2").Value = "6.39"
$q3_21.Range("G2").NumberFormat = "@"
$q3_21.Range("G2").Value = "0.0895"
$q3_21.Range("H2").Value = 4

# "2021-Q1" sheet is left completely untouched - its data/styling does not change.
